$wb = $excel.ActiveWorkbook

# --- Update selection on existing "Exam C" sheet ---
$wsC = $wb.Worksheets.Item("Exam C")
$wsC.Range("C6").Select()

# --- Update selection on existing "Exam B 2" sheet ---
$wsB2 = $wb.Worksheets.Item("Exam B 2")
$wsB2.Range("C59").Select()

# --- Add new "Exam C 2" sheet after "Exam B 2" (last sheet) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Exam C 2"

# --- Header reference tables (rows 1-6, columns D/E/G/I/K) ---
$ws.Range("D1").Value = "ALLOWED"
$ws.Range("E1").Value = "IPS"
$ws.Range("G1").Value = "DoS"
$ws.Range("I1").Value = "Diffusion"
$ws.Range("K1").Value = "RA"
$ws.Range("D2").Value = "BLOCKED"
$ws.Range("E2").Value = "Proxy"
$ws.Range("G2").Value = "Spoofing"
$ws.Range("I2").Value = "Confusion"
$ws.Range("K2").Value = "CRL"
$ws.Range("D3").Value = "ALLOWED"
$ws.Range("E3").Value = "Router"
$ws.Range("G3").Value = "Rootkit"
$ws.Range("I3").Value = "Collision"
$ws.Range("K3").Value = "CSR"
$ws.Range("D4").Value = "ALLOWED"
$ws.Range("E4").Value = "Load balancer"
$ws.Range("G4").Value = "Dictionary"
$ws.Range("I4").Value = "Obfuscation"
$ws.Range("K4").Value = "CA"
$ws.Range("D5").Value = "BLOCKED"
$ws.Range("E5").Value = "WAF"
$ws.Range("G5").Value = "Phishing"
$ws.Range("I5").Value = "Asymmetric"
$ws.Range("K5").Value = "OCSP"
$ws.Range("D6").Value = "ALLOWED"
$ws.Range("G6").Value = "Tailgating"
$ws.Range("I6").Value = "PFS"

# --- Score summary cells ---
$ws.Range("E11").Value = "85/85"
$ws.Range("E12").Value = "90/90"

# --- Answer grid: submitted answer (A), correct answer (B), grading formula (C) ---
$ws.Range("A6").Value = "D"
$ws.Range("B6").Value = "D"
$ws.Range("C6").Formula = '=IF(A6=B6, "Correct", "Incorrect")'
$ws.Range("A7").Value = "B"
$ws.Range("B7").Value = "B"
$ws.Range("C7").Formula = '=IF(A7=B7, "Correct", "Incorrect")'
$ws.Range("A8").Value = "D"
$ws.Range("B8").Value = "D"
$ws.Range("C8").Formula = '=IF(A8=B8, "Correct", "Incorrect")'
$ws.Range("A9").Value = "C"
$ws.Range("B9").Value = "C"
$ws.Range("C9").Formula = '=IF(A9=B9, "Correct", "Incorrect")'
$ws.Range("A10").Value = "D"
$ws.Range("B10").Value = "D"
$ws.Range("C10").Formula = '=IF(A10=B10, "Correct", "Incorrect")'
$ws.Range("A11").Value = "B"
$ws.Range("B11").Value = "B"
$ws.Range("C11").Formula = '=IF(A11=B11, "Correct", "Incorrect")'
$ws.Range("A12").Value = "B"
$ws.Range("B12").Value = "B"
$ws.Range("C12").Formula = '=IF(A12=B12, "Correct", "Incorrect")'
$ws.Range("A13").Value = "D"
$ws.Range("B13").Value = "D"
$ws.Range("C13").Formula = '=IF(A13=B13, "Correct", "Incorrect")'
$ws.Range("A14").Value = "D"
$ws.Range("B14").Value = "D"
$ws.Range("C14").Formula = '=IF(A14=B14, "Correct", "Incorrect")'
$ws.Range("A15").Value = "B"
$ws.Range("B15").Value = "B"
$ws.Range("C15").Formula = '=IF(A15=B15, "Correct", "Incorrect")'
$ws.Range("A16").Value = "CD"
$ws.Range("B16").Value = "CD"
$ws.Range("C16").Formula = '=IF(A16=B16, "Correct", "Incorrect")'
$ws.Range("A17").Value = "D"
$ws.Range("B17").Value = "D"
$ws.Range("C17").Formula = '=IF(A17=B17, "Correct", "Incorrect")'
$ws.Range("A18").Value = "A"
$ws.Range("B18").Value = "A"
$ws.Range("C18").Formula = '=IF(A18=B18, "Correct", "Incorrect")'
$ws.Range("A19").Value = "A"
$ws.Range("B19").Value = "A"
$ws.Range("C19").Formula = '=IF(A19=B19, "Correct", "Incorrect")'
$ws.Range("A20").Value = "A"
$ws.Range("B20").Value = "A"
$ws.Range("C20").Formula = '=IF(A20=B20, "Correct", "Incorrect")'
$ws.Range("A21").Value = "D"
$ws.Range("B21").Value = "D"
$ws.Range("C21").Formula = '=IF(A21=B21, "Correct", "Incorrect")'
$ws.Range("A22").Value = "C"
$ws.Range("B22").Value = "C"
$ws.Range("C22").Formula = '=IF(A22=B22, "Correct", "Incorrect")'
$ws.Range("A23").Value = "A"
$ws.Range("B23").Value = "A"
$ws.Range("C23").Formula = '=IF(A23=B23, "Correct", "Incorrect")'
$ws.Range("A24").Value = "C"
$ws.Range("B24").Value = "C"
$ws.Range("C24").Formula = '=IF(A24=B24, "Correct", "Incorrect")'
$ws.Range("A25").Value = "C"
$ws.Range("B25").Value = "C"
$ws.Range("C25").Formula = '=IF(A25=B25, "Correct", "Incorrect")'
$ws.Range("A26").Value = "A"
$ws.Range("B26").Value = "A"
$ws.Range("C26").Formula = '=IF(A26=B26, "Correct", "Incorrect")'
$ws.Range("A27").Value = "C"
$ws.Range("B27").Value = "C"
$ws.Range("C27").Formula = '=IF(A27=B27, "Correct", "Incorrect")'
$ws.Range("A28").Value = "B"
$ws.Range("B28").Value = "B"
$ws.Range("C28").Formula = '=IF(A28=B28, "Correct", "Incorrect")'
$ws.Range("A29").Value = "C"
$ws.Range("B29").Value = "C"
$ws.Range("C29").Formula = '=IF(A29=B29, "Correct", "Incorrect")'
$ws.Range("A30").Value = "B"
$ws.Range("B30").Value = "B"
$ws.Range("C30").Formula = '=IF(A30=B30, "Correct", "Incorrect")'
$ws.Range("A31").Value = "A"
$ws.Range("B31").Value = "A"
$ws.Range("C31").Formula = '=IF(A31=B31, "Correct", "Incorrect")'
$ws.Range("A32").Value = "B"
$ws.Range("B32").Value = "B"
$ws.Range("C32").Formula = '=IF(A32=B32, "Correct", "Incorrect")'
$ws.Range("A33").Value = "D"
$ws.Range("B33").Value = "D"
$ws.Range("C33").Formula = '=IF(A33=B33, "Correct", "Incorrect")'
$ws.Range("A34").Value = "D"
$ws.Range("B34").Value = "D"
$ws.Range("C34").Formula = '=IF(A34=B34, "Correct", "Incorrect")'
$ws.Range("A35").Value = "D"
$ws.Range("B35").Value = "D"
$ws.Range("C35").Formula = '=IF(A35=B35, "Correct", "Incorrect")'
$ws.Range("A36").Value = "C"
$ws.Range("B36").Value = "C"
$ws.Range("C36").Formula = '=IF(A36=B36, "Correct", "Incorrect")'
$ws.Range("A37").Value = "C"
$ws.Range("B37").Value = "C"
$ws.Range("C37").Formula = '=IF(A37=B37, "Correct", "Incorrect")'
$ws.Range("A38").Value = "B"
$ws.Range("B38").Value = "B"
$ws.Range("C38").Formula = '=IF(A38=B38, "Correct", "Incorrect")'
$ws.Range("A39").Value = "B"
$ws.Range("B39").Value = "B"
$ws.Range("C39").Formula = '=IF(A39=B39, "Correct", "Incorrect")'
$ws.Range("A40").Value = "B"
$ws.Range("B40").Value = "B"
$ws.Range("C40").Formula = '=IF(A40=B40, "Correct", "Incorrect")'
$ws.Range("A41").Value = "A"
$ws.Range("B41").Value = "A"
$ws.Range("C41").Formula = '=IF(A41=B41, "Correct", "Incorrect")'
$ws.Range("A42").Value = "B"
$ws.Range("B42").Value = "B"
$ws.Range("C42").Formula = '=IF(A42=B42, "Correct", "Incorrect")'
$ws.Range("A43").Value = "D"
$ws.Range("B43").Value = "D"
$ws.Range("C43").Formula = '=IF(A43=B43, "Correct", "Incorrect")'
$ws.Range("A44").Value = "A"
$ws.Range("B44").Value = "A"
$ws.Range("C44").Formula = '=IF(A44=B44, "Correct", "Incorrect")'
$ws.Range("A45").Value = "B"
$ws.Range("B45").Value = "B"
$ws.Range("C45").Formula = '=IF(A45=B45, "Correct", "Incorrect")'
$ws.Range("A46").Value = "D"
$ws.Range("B46").Value = "D"
$ws.Range("C46").Formula = '=IF(A46=B46, "Correct", "Incorrect")'
$ws.Range("A47").Value = "B"
$ws.Range("B47").Value = "B"
$ws.Range("C47").Formula = '=IF(A47=B47, "Correct", "Incorrect")'
$ws.Range("A48").Value = "C"
$ws.Range("B48").Value = "C"
$ws.Range("C48").Formula = '=IF(A48=B48, "Correct", "Incorrect")'
$ws.Range("A49").Value = "B"
$ws.Range("B49").Value = "B"
$ws.Range("C49").Formula = '=IF(A49=B49, "Correct", "Incorrect")'
$ws.Range("A50").Value = "C"
$ws.Range("B50").Value = "C"
$ws.Range("C50").Formula = '=IF(A50=B50, "Correct", "Incorrect")'
$ws.Range("A51").Value = "A"
$ws.Range("B51").Value = "A"
$ws.Range("C51").Formula = '=IF(A51=B51, "Correct", "Incorrect")'
$ws.Range("A52").Value = "A"
$ws.Range("B52").Value = "A"
$ws.Range("C52").Formula = '=IF(A52=B52, "Correct", "Incorrect")'
$ws.Range("A53").Value = "C"
$ws.Range("B53").Value = "C"
$ws.Range("C53").Formula = '=IF(A53=B53, "Correct", "Incorrect")'
$ws.Range("A54").Value = "DE"
$ws.Range("B54").Value = "DE"
$ws.Range("C54").Formula = '=IF(A54=B54, "Correct", "Incorrect")'
$ws.Range("A55").Value = "D"
$ws.Range("B55").Value = "D"
$ws.Range("C55").Formula = '=IF(A55=B55, "Correct", "Incorrect")'
$ws.Range("A56").Value = "B"
$ws.Range("B56").Value = "B"
$ws.Range("C56").Formula = '=IF(A56=B56, "Correct", "Incorrect")'
$ws.Range("A57").Value = "C"
$ws.Range("B57").Value = "C"
$ws.Range("C57").Formula = '=IF(A57=B57, "Correct", "Incorrect")'
$ws.Range("A58").Value = "A"
$ws.Range("B58").Value = "A"
$ws.Range("C58").Formula = '=IF(A58=B58, "Correct", "Incorrect")'
$ws.Range("A59").Value = "B"
$ws.Range("B59").Value = "B"
$ws.Range("C59").Formula = '=IF(A59=B59, "Correct", "Incorrect")'
$ws.Range("A60").Value = "E"
$ws.Range("B60").Value = "E"
$ws.Range("C60").Formula = '=IF(A60=B60, "Correct", "Incorrect")'
$ws.Range("A61").Value = "A"
$ws.Range("B61").Value = "A"
$ws.Range("C61").Formula = '=IF(A61=B61, "Correct", "Incorrect")'
$ws.Range("A62").Value = "A"
$ws.Range("B62").Value = "A"
$ws.Range("C62").Formula = '=IF(A62=B62, "Correct", "Incorrect")'
$ws.Range("A63").Value = "C"
$ws.Range("B63").Value = "C"
$ws.Range("C63").Formula = '=IF(A63=B63, "Correct", "Incorrect")'
$ws.Range("A64").Value = "B"
$ws.Range("B64").Value = "B"
$ws.Range("C64").Formula = '=IF(A64=B64, "Correct", "Incorrect")'
$ws.Range("A65").Value = "B"
$ws.Range("B65").Value = "B"
$ws.Range("C65").Formula = '=IF(A65=B65, "Correct", "Incorrect")'
$ws.Range("A66").Value = "A"
$ws.Range("B66").Value = "A"
$ws.Range("C66").Formula = '=IF(A66=B66, "Correct", "Incorrect")'
$ws.Range("A67").Value = "D"
$ws.Range("B67").Value = "D"
$ws.Range("C67").Formula = '=IF(A67=B67, "Correct", "Incorrect")'
$ws.Range("A68").Value = "A"
$ws.Range("B68").Value = "A"
$ws.Range("C68").Formula = '=IF(A68=B68, "Correct", "Incorrect")'
$ws.Range("A69").Value = "A"
$ws.Range("B69").Value = "A"
$ws.Range("C69").Formula = '=IF(A69=B69, "Correct", "Incorrect")'
$ws.Range("A70").Value = "B"
$ws.Range("B70").Value = "B"
$ws.Range("C70").Formula = '=IF(A70=B70, "Correct", "Incorrect")'
$ws.Range("A71").Value = "D"
$ws.Range("B71").Value = "D"
$ws.Range("C71").Formula = '=IF(A71=B71, "Correct", "Incorrect")'
$ws.Range("A72").Value = "B"
$ws.Range("B72").Value = "B"
$ws.Range("C72").Formula = '=IF(A72=B72, "Correct", "Incorrect")'
$ws.Range("A73").Value = "D"
$ws.Range("B73").Value = "D"
$ws.Range("C73").Formula = '=IF(A73=B73, "Correct", "Incorrect")'
$ws.Range("A74").Value = "C"
$ws.Range("B74").Value = "C"
$ws.Range("C74").Formula = '=IF(A74=B74, "Correct", "Incorrect")'
$ws.Range("A75").Value = "C"
$ws.Range("B75").Value = "C"
$ws.Range("C75").Formula = '=IF(A75=B75, "Correct", "Incorrect")'
$ws.Range("A76").Value = "C"
$ws.Range("B76").Value = "C"
$ws.Range("C76").Formula = '=IF(A76=B76, "Correct", "Incorrect")'
$ws.Range("A77").Value = "C"
$ws.Range("B77").Value = "C"
$ws.Range("C77").Formula = '=IF(A77=B77, "Correct", "Incorrect")'
$ws.Range("A78").Value = "A"
$ws.Range("B78").Value = "A"
$ws.Range("C78").Formula = '=IF(A78=B78, "Correct", "Incorrect")'
$ws.Range("A79").Value = "C"
$ws.Range("B79").Value = "C"
$ws.Range("C79").Formula = '=IF(A79=B79, "Correct", "Incorrect")'
$ws.Range("A80").Value = "C"
$ws.Range("B80").Value = "C"
$ws.Range("C80").Formula = '=IF(A80=B80, "Correct", "Incorrect")'
$ws.Range("A81").Value = "AD"
$ws.Range("B81").Value = "AD"
$ws.Range("C81").Formula = '=IF(A81=B81, "Correct", "Incorrect")'
$ws.Range("A82").Value = "BD"
$ws.Range("B82").Value = "BD"
$ws.Range("C82").Formula = '=IF(A82=B82, "Correct", "Incorrect")'
$ws.Range("A83").Value = "A"
$ws.Range("B83").Value = "A"
$ws.Range("C83").Formula = '=IF(A83=B83, "Correct", "Incorrect")'
$ws.Range("A84").Value = "D"
$ws.Range("B84").Value = "D"
$ws.Range("C84").Formula = '=IF(A84=B84, "Correct", "Incorrect")'
$ws.Range("A85").Value = "B"
$ws.Range("B85").Value = "B"
$ws.Range("C85").Formula = '=IF(A85=B85, "Correct", "Incorrect")'
$ws.Range("A86").Value = "A"
$ws.Range("B86").Value = "A"
$ws.Range("C86").Formula = '=IF(A86=B86, "Correct", "Incorrect")'
$ws.Range("A87").Value = "D"
$ws.Range("B87").Value = "D"
$ws.Range("C87").Formula = '=IF(A87=B87, "Correct", "Incorrect")'
$ws.Range("A88").Value = "C"
$ws.Range("B88").Value = "C"
$ws.Range("C88").Formula = '=IF(A88=B88, "Correct", "Incorrect")'
$ws.Range("A89").Value = "A"
$ws.Range("B89").Value = "A"
$ws.Range("C89").Formula = '=IF(A89=B89, "Correct", "Incorrect")'
$ws.Range("A90").Value = "C"
$ws.Range("B90").Value = "C"
$ws.Range("C90").Formula = '=IF(A90=B90, "Correct", "Incorrect")'

# --- Final selection on the new sheet ---
$ws.Range("H25").Select()
